# Update the weekly progress tracker sheet:
#  - E2: the "Trump API scraping" status note is revised to mention Tweepy
#  - A new status row (row 10) is appended for 羅佳敏's Mongo import / report work
#  - Selection ends on C10, matching where the new row's data was typed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Revise the note in E2
$ws.Range("E2").Value = "已註冊API並利用Tweepy/Get Old Tweets內建套件完成抓取"

# 2) Append new row 10, matching the date-formatted style already used in
#    column A (copy format only from A9, then set the value)
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 43824
$ws.Range("B10").Value = "羅佳敏"
$ws.Range("C10").Value = "協助Mongodb資料匯入、統整書面報告"
$ws.Range("D10").Value = "協助將抓取的資料匯入Mongodb，並統整期末書面報告"
$ws.Range("E10").Value = "全部已完成"

# 3) Leave the selection on C10, scrolled back to the top-left, like the
#    saved workbook
$ws.Range("A1").Select() | Out-Null
$ws.Range("C10").Select() | Out-Null
